# Export with no is_pref and no lev distance
# Re-orders the existing id/speaker_variant rows (2-24), clears the
# is_prefered ("x") markers in column D, and appends a new row (25)
# for the "Sillanus" entry, growing the sheet dimension to A1:H25.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = "#lelius"
$ws.Range("C2").Value = "Lelius"
$ws.Range("D2").Value = ""

$ws.Range("B3").Value = "#2.-soldaet"
$ws.Range("C3").Value = "2. Soldaet"
$ws.Range("D3").Value = ""

$ws.Range("B4").Value = "#mando"
$ws.Range("C4").Value = "Mando"
$ws.Range("D4").Value = ""

$ws.Range("B5").Value = "#l.-martius"
$ws.Range("C5").Value = "L. Martius"
$ws.Range("D5").Value = ""

$ws.Range("B6").Value = "#den-raed"
$ws.Range("C6").Value = "Den Raed"
$ws.Range("D6").Value = ""

$ws.Range("B7").Value = "#voester"
$ws.Range("C7").Value = "Voester"
$ws.Range("D7").Value = ""

$ws.Range("B8").Value = "#scipio-leest"
$ws.Range("C8").Value = "Scipio leest"
$ws.Range("D8").Value = ""

$ws.Range("B9").Value = "#sophon"
$ws.Range("C9").Value = "Sophon"
$ws.Range("D9").Value = ""

$ws.Range("B10").Value = "#m.-lelius"
$ws.Range("C10").Value = "M. Lelius"
$ws.Range("D10").Value = ""

$ws.Range("B11").Value = "#hasdru"
$ws.Range("C11").Value = "Hasdru"
$ws.Range("D11").Value = ""

$ws.Range("B12").Value = "#m.-leli"
$ws.Range("C12").Value = "M. Leli"
$ws.Range("D12").Value = ""

$ws.Range("B13").Value = "#q.-fabi"
$ws.Range("C13").Value = "Q. Fabi"
$ws.Range("D13").Value = ""

$ws.Range("B14").Value = "#scipio"
$ws.Range("C14").Value = "Scipio"
$ws.Range("D14").Value = ""

$ws.Range("B15").Value = "#masiniss"
$ws.Range("C15").Value = "Masiniss"
$ws.Range("D15").Value = ""

$ws.Range("B16").Value = "#2.sold"
$ws.Range("C16").Value = "2.Sold"
$ws.Range("D16").Value = ""

$ws.Range("B17").Value = "#bocher"
$ws.Range("C17").Value = "Bocher"
$ws.Range("D17").Value = ""

$ws.Range("B18").Value = "#methon"
$ws.Range("C18").Value = "Methon"
$ws.Range("D18").Value = ""

$ws.Range("B19").Value = "#metho"
$ws.Range("C19").Value = "Metho"
$ws.Range("D19").Value = ""

$ws.Range("B20").Value = "#1.-soldaet"
$ws.Range("C20").Value = "1. Soldaet"
$ws.Range("D20").Value = ""

$ws.Range("B21").Value = "#i"
$ws.Range("C21").Value = "I"
$ws.Range("D21").Value = ""

$ws.Range("B22").Value = "#siphax"
$ws.Range("C22").Value = "Siphax"
$ws.Range("D22").Value = ""

$ws.Range("B23").Value = "#babactus"
$ws.Range("C23").Value = "Babactus"
$ws.Range("D23").Value = ""

$ws.Range("B24").Value = "#luci.-m"
$ws.Range("C24").Value = "Luci. M"
$ws.Range("D24").Value = ""

$ws.Range("B25").Value = "#sillanus"
$ws.Range("C25").Value = "Sillanus"
$ws.Range("D25").Value = ""
$ws.Range("A25").Value = "https://www.dbnl.org/tekst/nieu001soph03_01"
$ws.Range("E25").Value = ""
$ws.Range("F25").Value = ""
$ws.Range("G25").Value = ""
$ws.Range("H25").Value = ""
